$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert 6 new rows before row 22 (after the last "YERSON" data row, which
#    is row 21). This pushes the existing row 22 (period 2112 for YERSON) down
#    to row 28, and the footer block (rows 27-28) down to rows 33-34.
# ---------------------------------------------------------------------------
$ws.Range("B22:B27").EntireRow.Insert()

# The new rows don't reliably inherit row 21's cell formatting (borders,
# number formats) from Insert() alone, so copy the formatting explicitly from
# the template data row (21) onto the freshly inserted rows (22-27).
$ws.Range("B21:J21").Copy()
$ws.Range("B22:J27").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Rows 16-21 now describe a NEW worker (NELSON ENRIQUE BALDIRIS LUNA),
#    replacing the data that used to belong to YERSON in that block.
# ---------------------------------------------------------------------------
$nelsonDoc = "1050944607"
$nelsonName = "NELSON ENRIQUE BALDIRIS LUNA"
$nelsonSalario = 828116

$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = $nelsonDoc
$ws.Range("D16").Value = $nelsonName
$ws.Range("E16").Value = "2112"
$ws.Range("F16").Value = 18726
$ws.Range("G16").Value = $nelsonSalario

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = $nelsonDoc
$ws.Range("D17").Value = $nelsonName
$ws.Range("E17").Value = "2111"
$ws.Range("F17").Value = 35112
$ws.Range("G17").Value = $nelsonSalario

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = $nelsonDoc
$ws.Range("D18").Value = $nelsonName
$ws.Range("E18").Value = "2110"
$ws.Range("F18").Value = 35112
$ws.Range("G18").Value = $nelsonSalario

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = $nelsonDoc
$ws.Range("D19").Value = $nelsonName
$ws.Range("E19").Value = "2109"
$ws.Range("F19").Value = 35112
$ws.Range("G19").Value = $nelsonSalario

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = $nelsonDoc
$ws.Range("D20").Value = $nelsonName
$ws.Range("E20").Value = "2108"
$ws.Range("F20").Value = 35112
$ws.Range("G20").Value = $nelsonSalario

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = $nelsonDoc
$ws.Range("D21").Value = $nelsonName
$ws.Range("E21").Value = "2107"
$ws.Range("F21").Value = 35112
$ws.Range("G21").Value = $nelsonSalario

# ---------------------------------------------------------------------------
# 3. Rows 22-28 keep describing YERSON ENRIQUE BALDIRIS LUNA, but his
#    "Salario Basico" (column G) is refreshed to the new value on every row.
# ---------------------------------------------------------------------------
$yersonDoc = "1047468022"
$yersonName = "YERSON ENRIQUE BALDIRIS LUNA"
$yersonSalario = 908526

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = $yersonDoc
$ws.Range("D22").Value = $yersonName
$ws.Range("E22").Value = "2112"
$ws.Range("F22").Value = 18726
$ws.Range("G22").Value = $yersonSalario

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = $yersonDoc
$ws.Range("D23").Value = $yersonName
$ws.Range("E23").Value = "2111"
$ws.Range("F23").Value = 35112
$ws.Range("G23").Value = $yersonSalario

$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = $yersonDoc
$ws.Range("D24").Value = $yersonName
$ws.Range("E24").Value = "2110"
$ws.Range("F24").Value = 36341
$ws.Range("G24").Value = $yersonSalario

$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = $yersonDoc
$ws.Range("D25").Value = $yersonName
$ws.Range("E25").Value = "2109"
$ws.Range("F25").Value = 36341
$ws.Range("G25").Value = $yersonSalario

$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = $yersonDoc
$ws.Range("D26").Value = $yersonName
$ws.Range("E26").Value = "2108"
$ws.Range("F26").Value = 36341
$ws.Range("G26").Value = $yersonSalario

$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = $yersonDoc
$ws.Range("D27").Value = $yersonName
$ws.Range("E27").Value = "2107"
$ws.Range("F27").Value = 36341
$ws.Range("G27").Value = $yersonSalario

$ws.Range("B28").Value = "CC"
$ws.Range("C28").Value = $yersonDoc
$ws.Range("D28").Value = $yersonName
$ws.Range("E28").Value = "2005"
$ws.Range("F28").Value = 15215
$ws.Range("G28").Value = $yersonSalario

# ---------------------------------------------------------------------------
# 4. Update the summary cells: total "VALOR MORA" and worker count.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 408703
$ws.Range("C13").Value = 2
